$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 26, shifting existing rows 26-29 down to 27-30.
$ws.Rows.Item(26).Insert()

# Match the date cell's number format (datetime format) used by the rest of column D.
$ws.Cells.Item(26, 4).NumberFormat = $ws.Cells.Item(27, 4).NumberFormat

# Populate the new row 26 with the new price record.
$ws.Cells.Item(26, 1).Value = 3
$ws.Cells.Item(26, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(26, 3).Value = "Coquimbo"
$ws.Cells.Item(26, 4).Value = 45015
$ws.Cells.Item(26, 5).Value = 5
$ws.Cells.Item(26, 6).Value = "Fruta"
$ws.Cells.Item(26, 7).Value = 100104
$ws.Cells.Item(26, 8).Value = "Frutos de pepita"
$ws.Cells.Item(26, 9).Value = 100104001
$ws.Cells.Item(26, 10).Value = "Granada"
$ws.Cells.Item(26, 11).Value = "Wonderfull"
$ws.Cells.Item(26, 12).Value = "Primera"
$ws.Cells.Item(26, 13).Value = 56
$ws.Cells.Item(26, 14).Value = 15000
$ws.Cells.Item(26, 15).Value = 15000
$ws.Cells.Item(26, 16).Value = 15000
$ws.Cells.Item(26, 17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(26, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(26, 19).Value = 1071
$ws.Cells.Item(26, 20).Value = 14
